$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Mustermann"
$ws.Range("C2").Value = "Midnight Rain"
$ws.Range("D2").Value = "Fantasy"
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "2000-12-16"
$ws.Range("F2").ClearFormats()
$ws.Range("G2").Value = "A former architect battles corporate zombies, 
      an evil sorceress, and her own childhood to become queen 
      of the world."

$ws.Range("C3").Value = "Maeve Ascendant"
$ws.Range("D3").Value = "Fantasy"
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "2000-11-17"
$ws.Range("F3").ClearFormats()
$ws.Range("G3").Value = "After the collapse of a nanotechnology 
      society in England, the young survivors lay the 
      foundation for a new society."
$ws.Range("H3").Value = "Corets, Eva"

$ws.Range("C4").Value = "Oberon's Legacy"
$ws.Range("D4").Value = "Fantasy"
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = "2001-03-10"
$ws.Range("F4").ClearFormats()
$ws.Range("G4").Value = "In post-apocalypse England, the mysterious 
      agent known only as Oberon helps to create a new life 
      for the inhabitants of London. Sequel to Maeve 
      Ascendant."
$ws.Range("H4").Value = "Corets, Eva"

$ws.Range("C5").Value = "The Sundered Grail"
$ws.Range("D5").Value = "Fantasy"
$ws.Range("G5").Value = "The two daughters of Maeve, half-sisters, 
      battle one another for control of England. Sequel to 
      Oberon's Legacy."
$ws.Range("H5").Value = "Corets, Eva"

$ws.Range("C6").Value = "Lover Birds"
$ws.Range("D6").Value = "Romance"
$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").Value = "2000-09-02"
$ws.Range("F6").ClearFormats()
$ws.Range("G6").Value = "When Carla meets Paul at an ornithology 
      conference, tempers fly as feathers get ruffled."
$ws.Range("H6").Value = "Randall, Cynthia"

$ws.Range("C7").Value = "Splish Splash"
$ws.Range("D7").Value = "Romance"
$ws.Range("F7").NumberFormat = "@"
$ws.Range("F7").Value = "2000-11-02"
$ws.Range("F7").ClearFormats()
$ws.Range("G7").Value = "A deep sea diver finds true love twenty 
      thousand leagues beneath the sea."
$ws.Range("H7").Value = "Thurman, Paula"

$ws.Range("C8").Value = "Creepy Crawlies"
$ws.Range("D8").Value = "Horror"
$ws.Range("F8").NumberFormat = "@"
$ws.Range("F8").Value = "2000-12-06"
$ws.Range("F8").ClearFormats()
$ws.Range("G8").Value = "An anthology of horror stories about roaches,
      centipedes, scorpions  and other insects."
$ws.Range("H8").Value = "Knorr, Stefan"

$ws.Range("C9").Value = "Paradox Lost"
$ws.Range("D9").Value = "Science Fiction"
$ws.Range("F9").NumberFormat = "@"
$ws.Range("F9").Value = "2000-11-02"
$ws.Range("F9").ClearFormats()
$ws.Range("G9").Value = "After an inadvertant trip through a Heisenberg
      Uncertainty Device, James Salway discovers the problems 
      of being quantum."
$ws.Range("H9").Value = "Kress, Peter"

$ws.Range("C10").Value = "Microsoft .NET: The Programming Bible"
$ws.Range("D10").Value = "Computer"
$ws.Range("F10").NumberFormat = "@"
$ws.Range("F10").Value = "2000-12-09"
$ws.Range("F10").ClearFormats()
$ws.Range("G10").Value = "Microsoft's .NET initiative is explored in 
      detail in this deep programmer's reference."
$ws.Range("H10").Value = "O'Brien, Tim"

$ws.Range("C11").Value = "MSXML3: A Comprehensive Guide"
$ws.Range("D11").Value = "Computer"
$ws.Range("F11").NumberFormat = "@"
$ws.Range("F11").Value = "2000-12-01"
$ws.Range("F11").ClearFormats()
$ws.Range("G11").Value = "The Microsoft MSXML3 parser is covered in 
      detail, with attention to XML DOM interfaces, XSLT processing, 
      SAX and more."
$ws.Range("H11").Value = "O'Brien, Tim"

$ws.Range("C12").Value = "Visual Studio 7: A Comprehensive Guide"
$ws.Range("D12").Value = "Computer"
$ws.Range("F12").NumberFormat = "@"
$ws.Range("F12").Value = "2001-04-16"
$ws.Range("F12").ClearFormats()
$ws.Range("G12").Value = "Microsoft Visual Studio 7 is explored in depth,
      looking at how Visual Basic, Visual C++, C#, and ASP+ are 
      integrated into a comprehensive development 
      environment."
$ws.Range("H12").Value = "Galos, Mike"

